$wb = $excel.ActiveWorkbook

$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A6").NumberFormat = "@"
$wsGaz.Range("A6").Value = "2025-06-20"
$wsGaz.Range("A6").ClearFormats()
$wsGaz.Range("B6").Value = 39.7

$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A6").NumberFormat = "@"
$wsCO2.Range("A6").Value = "2025-06-20"
$wsCO2.Range("A6").ClearFormats()
$wsCO2.Range("B6").Value = 72.2
